# Updates symbol list data (prices / 1h volume % / a few coin rows that got
# reordered) to match the "Updated symbol list" GitHub Actions commit.
#
# Columns D (Price) and E (Volume(1h)) hold numeric-/percent-looking text
# that must remain plain text (as it was originally stored, e.g. "247.19",
# "0.81%"). A leading apostrophe forces Excel to keep the entry as text
# instead of auto-converting it to a number/percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.28"
$ws.Range("E2").Value = "'0.89%"
$ws.Range("E3").Value = "'4.28%"
$ws.Range("D5").Value = "'0.05601"
$ws.Range("E5").Value = "'-0.34%"
$ws.Range("D6").Value = "'6.474"
$ws.Range("E6").Value = "'-1.48%"
$ws.Range("D7").Value = "'0.8130"
$ws.Range("E7").Value = "'-0.11%"
$ws.Range("D8").Value = "'0.8439"
$ws.Range("E8").Value = "'0.63%"
$ws.Range("D9").Value = "'0.06989"
$ws.Range("E9").Value = "'0.57%"
$ws.Range("E10").Value = "'-1.07%"
$ws.Range("D11").Value = "'0.09396"
$ws.Range("E11").Value = "'-0.05%"
$ws.Range("E12").Value = "'0.58%"
$ws.Range("D13").Value = "'0.0005970"
$ws.Range("E13").Value = "'-93.86%"
$ws.Range("D14").Value = "'0.006188"
$ws.Range("E14").Value = "'1.28%"
$ws.Range("D15").Value = "'3.606"
$ws.Range("E15").Value = "'3.06%"
$ws.Range("D16").Value = "'3.018"
$ws.Range("E16").Value = "'0.32%"
$ws.Range("E17").Value = "'-1.73%"
$ws.Range("D18").Value = "'0.3121"
$ws.Range("E18").Value = "'-1.95%"
$ws.Range("D19").Value = "'0.1338"
$ws.Range("E19").Value = "'0.04%"
$ws.Range("D20").Value = "'0.03209"
$ws.Range("E20").Value = "'-1.69%"
$ws.Range("E21").Value = "'-1.25%"
$ws.Range("D22").Value = "'3.744"
$ws.Range("E22").Value = "'0.02%"
$ws.Range("D23").Value = "'0.04690"
$ws.Range("E23").Value = "'0.62%"
$ws.Range("E24").Value = "'-1.34%"
$ws.Range("E25").Value = "'0.00%"
$ws.Range("D26").Value = "'0.004574"
$ws.Range("E26").Value = "'0.92%"
$ws.Range("D27").Value = "'0.00009598"
$ws.Range("E27").Value = "'-0.94%"
$ws.Range("E28").Value = "'0.01%"
$ws.Range("D40").Value = "'0.03662"
$ws.Range("E40").Value = "'-0.05%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1057"
$ws.Range("E41").Value = "'0.15%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002499"
$ws.Range("E42").Value = "'-8.51%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.006149"
$ws.Range("E43").Value = "'-1.04%"
$ws.Range("D44").Value = "'0.008287"
$ws.Range("E44").Value = "'1.43%"
$ws.Range("D45").Value = "'0.00005401"
$ws.Range("E45").Value = "'2.13%"
$ws.Range("E46").Value = "'0.09%"
$ws.Range("E47").Value = "'-35.77%"
$ws.Range("D48").Value = "'0.002584"
$ws.Range("E48").Value = "'27.90%"
$ws.Range("E49").Value = "'0.09%"
$ws.Range("E50").Value = "'0.09%"
